# Applies the "014 och 096 klar" update to Rawdata/lydata.xlsx
#  - bumps the workbook window x-position
#  - appends 10 measurement rows for lya zz014 and 10 for lya zz096
#  - updates the frozen-pane/selection to reflect the new bottom row
#  - dimension / sharedStrings bookkeeping is handled automatically by Excel

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- workbook window position -------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Left = 2800

# --- comment text for the two new lyor -----------------------------------------
$comment014 = "Om man räknar med nysnö är andel snöfri yta 70 %. Den snöfria ytan är då 236,6 kvadratmeter. Det fanns 11 hål på lyan men alla utom ett var täckt med snö. Jag borde ha undersökt om det bara var nysnö men det gjorde jag inte."
$comment096 = 'Lyan ligger i en sluttning med riktning 172 grader S. Dock skymmer en liten kulle lyan. Kullen har inga synliga öppningar och väldigt lite högt gräs. Lyans riktning blir snarare 266 grader. Vinkel svårt pga snö. Kullen framför lyan har en vinkel på 28 grader. Det är ett mycket färskt rävspår på lyan och ett grävt hål i snön bakom en sten 10 m norr om "kullen". Endast 1 synlig öppning som är igenisad. Andel snöfri yta är svårt att säga. Kullens snöfria area är 18,7 kvadratmeter men oklart om den är en del av lyan.'

# --- new data rows ---------------------------------------------------------------
# columns: A lya, B snöfri procent, C snöfri area, D lufttemp,
#          E marktemp orange, F marktemp svart, G underlag marktemp, H snödjup,
#          I riktning, J vinkel, K antal lyöppningar, L aktiv, M kommentar
$rows = @(
    @{A="zz014"; B=90; C=449.5; D=1.2; E=0.9;  F=0.5;  G="b";  H=0;  I="130 SO"; J=9;  K=1; L="j"; M=$comment014}
    @{A="zz014"; B=90; C=449.5; D=1.2; E=-0.1; F=0.1;  G="ns"; H=5;  I="130 SO"; J=9;  K=1; L="j"; M=$comment014}
    @{A="zz014"; B=90; C=449.5; D=1.2; E=-0.3; F=-0.4; G="s";  H=12; I="130 SO"; J=9;  K=1; L="j"; M=$comment014}
    @{A="zz014"; B=90; C=449.5; D=1.2; E=-0.2; F=-0.3; G="ns"; H=12; I="130 SO"; J=9;  K=1; L="j"; M=$comment014}
    @{A="zz014"; B=90; C=449.5; D=1.2; E=-0.2; F=-0.4; G="b";  H=0;  I="130 SO"; J=9;  K=1; L="j"; M=$comment014}
    @{A="zz014"; B=90; C=449.5; D=1.2; E=4.6;  F=4.5;  G="b";  H=0;  I="130 SO"; J=9;  K=1; L="j"; M=$comment014}
    @{A="zz014"; B=90; C=449.5; D=1.2; E=0.5;  F=-0.3; G="b";  H=0;  I="130 SO"; J=9;  K=1; L="j"; M=$comment014}
    @{A="zz014"; B=90; C=449.5; D=1.2; E=-0.3; F=-0.4; G="ns"; H=3;  I="130 SO"; J=9;  K=1; L="j"; M=$comment014}
    @{A="zz014"; B=90; C=449.5; D=1.2; E=-0.3; F=-0.3; G="ns"; H=3;  I="130 SO"; J=9;  K=1; L="j"; M=$comment014}
    @{A="zz014"; B=90; C=449.5; D=1.2; E=-0.3; F=-0.3; G="ns"; H=10; I="130 SO"; J=9;  K=1; L="j"; M=$comment014}

    @{A="zz096"; B=15; C=10.3;  D=0.1; E=-0.3; F=-0.4; G="s";  H=59; I="266 V";  J=16; K=0; L="n"; M=$comment096}
    @{A="zz096"; B=15; C=10.3;  D=0.1; E=-0.3; F=-0.4; G="s";  H=16; I="266 V";  J=16; K=0; L="n"; M=$comment096}
    @{A="zz096"; B=15; C=10.3;  D=0.1; E=-0.3; F=-0.3; G="s";  H=14; I="266 V";  J=16; K=0; L="n"; M=$comment096}
    @{A="zz096"; B=15; C=10.3;  D=0.1; E=-0.5; F=-0.5; G="b";  H=0;  I="266 V";  J=16; K=0; L="n"; M=$comment096}
    @{A="zz096"; B=15; C=10.3;  D=0.1; E=-0.4; F=-0.5; G="s";  H=28; I="266 V";  J=16; K=0; L="n"; M=$comment096}
    @{A="zz096"; B=15; C=10.3;  D=0.1; E=-0.2; F=-0.3; G="b";  H=0;  I="266 V";  J=16; K=0; L="n"; M=$comment096}
    @{A="zz096"; B=15; C=10.3;  D=0.1; E=-0.3; F=-0.4; G="b";  H=0;  I="266 V";  J=16; K=0; L="n"; M=$comment096}
    @{A="zz096"; B=15; C=10.3;  D=0.1; E=-0.3; F=-0.4; G="b";  H=0;  I="266 V";  J=16; K=0; L="n"; M=$comment096}
    @{A="zz096"; B=15; C=10.3;  D=0.1; E=-0.2; F=-0.3; G="b";  H=0;  I="266 V";  J=16; K=0; L="n"; M=$comment096}
    @{A="zz096"; B=15; C=10.3;  D=0.1; E=-0.3; F=-0.4; G="s";  H=29; I="266 V";  J=16; K=0; L="n"; M=$comment096}
)

$startRow = 82
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $d = $rows[$i]
    $ws.Cells.Item($r, 1).Value  = $d.A
    $ws.Cells.Item($r, 2).Value  = $d.B
    $ws.Cells.Item($r, 3).Value  = $d.C
    $ws.Cells.Item($r, 4).Value  = $d.D
    $ws.Cells.Item($r, 5).Value  = $d.E
    $ws.Cells.Item($r, 6).Value  = $d.F
    $ws.Cells.Item($r, 7).Value  = $d.G
    $ws.Cells.Item($r, 8).Value  = $d.H
    $ws.Cells.Item($r, 9).Value  = $d.I
    $ws.Cells.Item($r, 10).Value = $d.J
    $ws.Cells.Item($r, 11).Value = $d.K
    $ws.Cells.Item($r, 12).Value = $d.L
    $ws.Cells.Item($r, 13).Value = $d.M
}

$lastRow = $startRow + $rows.Count - 1

# --- refresh the view: keep header row frozen, scroll/select the new bottom row --
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("G" + $lastRow).Select()
